$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value forcing text interpretation (no auto numeric coercion),
# then reset style back to Normal so no stray number-format/style is introduced.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.323.83'
$ws.Range("E2").Value = "  -3.43%  "
Set-TextValue 'D3' '1.933.92'
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue 'D5' '249.30'
$ws.Range("E5").Value = "  -4.16%  "
Set-TextValue 'D6' '0.7238'
$ws.Range("E6").Value = "  -8.39%  "
Set-TextValue 'D7' '1.000'
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue 'D8' '0.3296'
$ws.Range("E8").Value = "  -8.75%  "
Set-TextValue 'D9' '27.72'
$ws.Range("E9").Value = "  -3.42%  "
Set-TextValue 'D10' '0.06846'
$ws.Range("E10").Value = "  -3.34%  "
Set-TextValue 'D11' '0.8063'
$ws.Range("E11").Value = "  -5.40%  "
Set-TextValue 'D12' '0.08074'
$ws.Range("E12").Value = "  -0.38%  "
Set-TextValue 'D13' '1.932.93'
$ws.Range("E13").Value = "  -3.32%  "
Set-TextValue 'D14' '5.421'
$ws.Range("E14").Value = "  -3.83%  "
Set-TextValue 'D15' '94.96'
$ws.Range("E15").Value = "  -6.57%  "
Set-TextValue 'D16' '14.53'
$ws.Range("E16").Value = "  -1.28%  "
Set-TextValue 'D17' '30.329.87'
$ws.Range("E17").Value = "  -3.39%  "
Set-TextValue 'D18' '251.93'
$ws.Range("E18").Value = "  -8.78%  "
Set-TextValue 'D19' '0.000008050'
$ws.Range("E19").Value = "  +1.36%  "
Set-TextValue 'D20' '5.826'
$ws.Range("E20").Value = "  -2.06%  "
Set-TextValue 'D21' '2.188.08'
$ws.Range("E21").Value = "  -3.36%  "
Set-TextValue 'D22' '1.000'
$ws.Range("E22").Value = "  +0.12%  "
Set-TextValue 'D23' '1.000'
$ws.Range("E23").Value = "  +0.05%  "
Set-TextValue 'D24' '6.879'
$ws.Range("E24").Value = "  -4.72%  "
Set-TextValue 'D25' '9.716'
$ws.Range("E25").Value = "  -4.85%  "
Set-TextValue 'D26' '159.70'
$ws.Range("E26").Value = "  -3.10%  "
Set-TextValue 'D27' '2.392'
$ws.Range("E27").Value = "  -0.31%  "
Set-TextValue 'D28' '19.14'
$ws.Range("E28").Value = "  -4.91%  "
Set-TextValue 'D29' '0.1343'
$ws.Range("E29").Value = "  -10.94%  "
Set-TextValue 'D30' '1.557'
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("E31").Value = "  -1.68%  "
Set-TextValue 'D32' '4.402'
$ws.Range("E32").Value = "  -5.25%  "
Set-TextValue 'D33' '4.175'
$ws.Range("E33").Value = "  -5.38%  "
Set-TextValue 'D34' '0.05099'
$ws.Range("E34").Value = "  -2.57%  "
Set-TextValue 'D35' '1.222'
$ws.Range("E35").Value = "  -0.42%  "
Set-TextValue 'D36' '0.7411'
$ws.Range("E36").Value = "  -3.77%  "
Set-TextValue 'D37' '2.751'
$ws.Range("E37").Value = "  -2.15%  "
Set-TextValue 'D38' '0.01970'
$ws.Range("E38").Value = "  -2.20%  "
Set-TextValue 'D39' '2.836'
$ws.Range("E39").Value = "  -4.05%  "
Set-TextValue 'D40' '6.598'
$ws.Range("E40").Value = "  -1.77%  "
Set-TextValue 'D41' '79.24'
$ws.Range("E41").Value = "  -2.72%  "
Set-TextValue 'D42' '0.4470'
$ws.Range("E42").Value = "  -5.96%  "
Set-TextValue 'D43' '1.996'
$ws.Range("E43").Value = "  -9.23%  "
$ws.Range("E44").Value = "  +0.03%  "
Set-TextValue 'D45' '0.8351'
$ws.Range("E45").Value = "  -2.94%  "
Set-TextValue 'D46' '102.00'
$ws.Range("E46").Value = "  -2.87%  "
Set-TextValue 'D47' '9.771'
$ws.Range("E47").Value = "  -1.99%  "
Set-TextValue 'D48' '7.303'
$ws.Range("E48").Value = "  -4.99%  "
Set-TextValue 'D49' '36.52'
$ws.Range("E49").Value = "  -1.45%  "
Set-TextValue 'D50' '0.05958'
$ws.Range("E50").Value = "  -0.28%  "
Set-TextValue 'D51' '1.476'
$ws.Range("E51").Value = "  -1.03%  "
